$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.730.48"
$ws.Range("E2").Value = "  +2.52%  "
$ws.Range("D3").Value = "1.873.48"
$ws.Range("E3").Value = "  +2.28%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "324.50"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "0.4610"
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("D8").Value = "0.3871"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "0.07861"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").Value = "0.9882"
$ws.Range("E10").Value = "  +2.98%  "
$ws.Range("D11").Value = "21.79"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").Value = "1.915.33"
$ws.Range("E12").Value = "  +5.71%  "
$ws.Range("D13").Value = "6.986"
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("D14").Value = "5.711"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").Value = "0.06981"
$ws.Range("E15").Value = "  +1.80%  "
$ws.Range("D16").Value = "88.44"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "0.00001002"
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("D19").Value = "16.79"
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").Value = "28.739.57"
$ws.Range("E21").Value = "  +2.51%  "
$ws.Range("D22").Value = "5.282"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "11.04"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").Value = "2.101"
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").Value = "2.094.34"
$ws.Range("E25").Value = "  +3.04%  "
$ws.Range("D26").Value = "152.68"
$ws.Range("E26").Value = "  -1.42%  "
$ws.Range("D27").Value = "19.27"
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("D28").Value = "5.850"
$ws.Range("E28").Value = "  +3.55%  "
$ws.Range("D29").Value = "1.976"
$ws.Range("E29").Value = "  +0.99%  "
$ws.Range("D30").Value = "118.87"
$ws.Range("E30").Value = "  +0.61%  "
$ws.Range("D31").Value = "0.09316"
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("D32").Value = "0.9207"
$ws.Range("E32").Value = "  -1.29%  "
$ws.Range("D33").Value = "5.305"
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("D34").Value = "1.338"
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").Value = "0.05775"
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("D37").Value = "1.153"
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("D38").Value = "0.02072"
$ws.Range("E38").Value = "  -2.35%  "
$ws.Range("D39").Value = "7.684"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").Value = "0.5642"
$ws.Range("E40").Value = "  +0.99%  "
$ws.Range("E41").Value = "  +1.61%  "
$ws.Range("D42").Value = "9.833"
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").Value = "0.07212"
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("D44").Value = "11.72"
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("D45").Value = "0.5294"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").Value = "2.130"
$ws.Range("E46").Value = "  +1.63%  "
$ws.Range("E47").Value = "  -2.47%  "
$ws.Range("D48").Value = "1.832"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("D49").Value = "113.48"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").Value = "2.418"
$ws.Range("E50").Value = "  +4.02%  "
$ws.Range("D51").Value = "1.003"
$ws.Range("E51").Value = "  +0.27%  "
